$d = $word.ActiveDocument

$replacements = @(
    @("741÷6=", "373÷6="),
    @("145÷6=", "540÷5="),
    @("744÷4=", "149÷2="),
    @("281÷4=", "745÷6="),
    @("205÷2=", "894÷2="),
    @("324÷5=", "455÷4="),
    @("973÷5=", "579÷8="),
    @("434÷7=", "373÷3="),
    @("294÷9=", "614÷6="),
    @("394÷3=", "310÷9="),
    @("372÷6=", "992÷6="),
    @("768÷7=", "710÷3="),
    @("477÷5=", "600÷4="),
    @("362÷8=", "906÷6="),
    @("304÷5=", "145÷7="),
    @("866÷4=", "898÷2="),
    @("423÷6=", "209÷7="),
    @("798÷2=", "571÷4="),
    @("138÷6=", "598÷3="),
    @("489÷6=", "208÷6="),
    @("597÷8=", "800÷6="),
    @("119÷3=", "270÷9="),
    @("851÷5=", "391÷5="),
    @("224÷2=", "501÷2="),
    @("887÷4=", "911÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
